# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns on Sheet1
# with newly scraped values (rows 2-51), as produced by the scheduled
# GitHub Actions job that regenerates cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to stay a text value (matches the source data which
    # always stores these columns as strings, e.g. "1.00", "64.842.94"),
    # then restore the default (Normal) style so no formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "64.842.94"
Set-TextValue $ws.Range("E2") "  +0.90%  "
# Row 3
Set-TextValue $ws.Range("D3") "3.145.88"
Set-TextValue $ws.Range("E3") "  +1.51%  "
# Row 4
Set-TextValue $ws.Range("E4") "  +0.07%  "
# Row 5
Set-TextValue $ws.Range("D5") "572.26"
Set-TextValue $ws.Range("E5") "  +2.20%  "
# Row 6
Set-TextValue $ws.Range("D6") "151.15"
Set-TextValue $ws.Range("E6") "  +4.42%  "
# Row 7
Set-TextValue $ws.Range("E7") "  +0.04%  "
# Row 8
Set-TextValue $ws.Range("D8") "3.142.01"
Set-TextValue $ws.Range("E8") "  +1.45%  "
# Row 9
Set-TextValue $ws.Range("E9") "  +4.19%  "
# Row 10
Set-TextValue $ws.Range("E10") "  +4.81%  "
# Row 11
Set-TextValue $ws.Range("D11") "6.18"
Set-TextValue $ws.Range("E11") "  +0.18%  "
# Row 12
Set-TextValue $ws.Range("E12") "  +6.40%  "
# Row 13
Set-TextValue $ws.Range("E13") "  +10.93%  "
# Row 14
Set-TextValue $ws.Range("D14") "37.53"
Set-TextValue $ws.Range("E14") "  +6.54%  "
# Row 15
Set-TextValue $ws.Range("D15") "3.662.24"
Set-TextValue $ws.Range("E15") "  +2.02%  "
# Row 16
Set-TextValue $ws.Range("D16") "64.914.31"
Set-TextValue $ws.Range("E16") "  +0.97%  "
# Row 17
Set-TextValue $ws.Range("E17") "  +6.75%  "
# Row 18
Set-TextValue $ws.Range("D18") "3.148.73"
Set-TextValue $ws.Range("E18") "  +1.87%  "
# Row 19
Set-TextValue $ws.Range("E19") "  +0.29%  "
# Row 20
Set-TextValue $ws.Range("D20") "512.08"
Set-TextValue $ws.Range("E20") "  +5.69%  "
# Row 21
Set-TextValue $ws.Range("D21") "14.92"
Set-TextValue $ws.Range("E21") "  +6.66%  "
# Row 22
Set-TextValue $ws.Range("D22") "0.735"
Set-TextValue $ws.Range("E22") "  +8.82%  "
# Row 23
Set-TextValue $ws.Range("D23") "15.32"
Set-TextValue $ws.Range("E23") "  +9.63%  "
# Row 24
Set-TextValue $ws.Range("E24") "  +3.52%  "
# Row 25
Set-TextValue $ws.Range("D25") "85.04"
Set-TextValue $ws.Range("E25") "  +4.56%  "
# Row 26
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.05%  "
# Row 27
Set-TextValue $ws.Range("E27") "  +4.21%  "
# Row 28
Set-TextValue $ws.Range("E28") "  +8.43%  "
# Row 29
Set-TextValue $ws.Range("E29") "  +5.31%  "
# Row 30
Set-TextValue $ws.Range("E30") "  +6.25%  "
# Row 31
Set-TextValue $ws.Range("D31") "1.00"
Set-TextValue $ws.Range("E31") "  +0.11%  "
# Row 32
Set-TextValue $ws.Range("E32") "  +3.22%  "
# Row 33
Set-TextValue $ws.Range("E33") "  +6.19%  "
# Row 34
Set-TextValue $ws.Range("D34") "6.09"
Set-TextValue $ws.Range("E34") "  +8.49%  "
# Row 35
Set-TextValue $ws.Range("E35") "  +5.77%  "
# Row 36
Set-TextValue $ws.Range("D36") "55.41"
Set-TextValue $ws.Range("E36") "  -0.89%  "
# Row 37
Set-TextValue $ws.Range("D37") "482.86"
Set-TextValue $ws.Range("E37") "  +5.40%  "
# Row 38
Set-TextValue $ws.Range("D38") "0.0862"
Set-TextValue $ws.Range("E38") "  +4.92%  "
# Row 39
Set-TextValue $ws.Range("D39") "0.0423"
Set-TextValue $ws.Range("E39") "  +3.63%  "
# Row 40
Set-TextValue $ws.Range("D40") "3.01"
Set-TextValue $ws.Range("E40") "  -0.47%  "
# Row 41
Set-TextValue $ws.Range("D41") "3.115.16"
Set-TextValue $ws.Range("E41") "  +4.32%  "
# Row 42
Set-TextValue $ws.Range("D42") "8.64"
Set-TextValue $ws.Range("E42") "  +4.49%  "
# Row 43
Set-TextValue $ws.Range("E43") "  +4.02%  "
# Row 44
Set-TextValue $ws.Range("D44") "0.292"
Set-TextValue $ws.Range("E44") "  +11.19%  "
# Row 45
Set-TextValue $ws.Range("E45") "  +15.19%  "
# Row 46
Set-TextValue $ws.Range("D46") "29.07"
Set-TextValue $ws.Range("E46") "  +3.80%  "
# Row 47
Set-TextValue $ws.Range("D47") "0.0₃0572"
Set-TextValue $ws.Range("E47") "  +10.74%  "
# Row 48
Set-TextValue $ws.Range("E48") "  +0.01%  "
# Row 49
Set-TextValue $ws.Range("E49") "  +3.07%  "
# Row 50
Set-TextValue $ws.Range("E50") "  +10.01%  "
# Row 51
Set-TextValue $ws.Range("D51") "118.66"
Set-TextValue $ws.Range("E51") "  -1.34%  "
